$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Runtimes")

$ws.Range("E4").Value = 0.0007986111111111112
$ws.Range("E5").Value = 0.00017361111111111112
$ws.Range("F5").Value = 0.0007986111111111112

$excel.Calculate()
